$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 226; this shifts rows 226:325 down to 227:326
$ws.Rows("226:226").Insert()

# Populate the newly inserted row 226 with its data (mirrors the
# surrounding/constant columns and sets the new record's values)
$ws.Range("A226").Value = 5
$ws.Range("B226").Value = "Macroferia Regional de Talca"
$ws.Range("C226").Value = "Maule"
$ws.Range("D226").Value = 44992
$ws.Range("D226").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E226").Value = 7
$ws.Range("F226").Value = 100112024
$ws.Range("G226").Value = "Choclo"
$ws.Range("H226").Value = "Choclero"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 20000
$ws.Range("K226").Value = 500
$ws.Range("L226").Value = 500
$ws.Range("M226").Value = 500
$ws.Range("N226").Value = "$/unidad"
$ws.Range("O226").Value = "Región del Maule"
$ws.Range("P226").Value = 500
$ws.Range("Q226").Value = 1
$ws.Range("R226").Value = "Hortaliza"
